$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('PayPal')
$ws.Cells.Item(2,1).Value = '''12/05/2023'
$ws.Cells.Item(2,2).Value = 'Payment Refund: Upwork Escrow Inc. USD 107.95'
$ws.Cells.Item(2,3).Value = '''107.95'
$ws.Cells.Item(2,4).Value = 'PayPal'
$ws.Cells.Item(3,1).Value = '''12/31/2023'
$ws.Cells.Item(3,2).Value = 'Cryptocurrency Transfer to 0x9C4d2...4f50 ETH -0.04347920'
$ws.Cells.Item(3,3).Value = '''-0.04371186'
$ws.Cells.Item(3,4).Value = 'PayPal'
$ws.Cells.Item(4,1).Value = '''12/05/2023'
$ws.Cells.Item(4,2).Value = 'User Initiated Withdrawal USD -107.95'
$ws.Cells.Item(4,3).Value = '''-107.95'
$ws.Cells.Item(4,4).Value = 'PayPal'
$ws.Cells.Item(5,1).Value = '''12/28/23'
$ws.Cells.Item(5,2).Value = 'Payment'
$ws.Cells.Item(5,3).Value = '''-85.57'
$ws.Cells.Item(5,4).Value = 'PayPal'
$ws.Cells.Item(6,1).Value = '''12/03/2023'
$ws.Cells.Item(6,2).Value = 'PreApproved Payment Bill User Payment: USD -52.50'
$ws.Cells.Item(6,3).Value = '''-52.50'
$ws.Cells.Item(6,4).Value = 'PayPal'
$ws.Cells.Item(7,1).Value = '''12/31/2023'
$ws.Cells.Item(7,2).Value = 'Cryptocurrency 0.04371186 ETH Buy USD -100.54'
$ws.Cells.Item(7,3).Value = '''-104.03'
$ws.Cells.Item(7,4).Value = 'PayPal'
$ws.Cells.Item(8,1).Value = '''12/05/2023'
$ws.Cells.Item(8,2).Value = 'PreApproved Payment Bill User Payment: USD -76.45'
$ws.Cells.Item(8,3).Value = '''-76.45'
$ws.Cells.Item(8,4).Value = 'PayPal'
$ws.Cells.Item(9,1).Value = '''12/28/2023'
$ws.Cells.Item(9,2).Value = 'PreApproved Payment Bill User Payment: USD -38.15'
$ws.Cells.Item(9,3).Value = '''-38.15'
$ws.Cells.Item(9,4).Value = 'PayPal'
$ws.Cells.Item(10,1).Value = '''12/04/2023'
$ws.Cells.Item(10,2).Value = 'PreApproved Payment Bill User Payment: USD -20.92'
$ws.Cells.Item(10,3).Value = '''-20.92'
$ws.Cells.Item(10,4).Value = 'PayPal'
$ws.Cells.Item(11,1).Value = '''12/19/2023'
$ws.Cells.Item(11,2).Value = 'PreApproved Payment Bill User Payment: USD -50.00'
$ws.Cells.Item(11,3).Value = '''-50.00'
$ws.Cells.Item(11,4).Value = 'PayPal'
$ws.Cells.Item(12,1).Value = '''12/28/2023'
$ws.Cells.Item(12,2).Value = 'PreApproved Payment Bill User Payment: USD -262.50'
$ws.Cells.Item(12,3).Value = '''-262.50'
$ws.Cells.Item(12,4).Value = 'PayPal'
$ws.Cells.Item(13,1).Value = '''12/30/2023'
$ws.Cells.Item(13,2).Value = 'PreApproved Payment Bill User Payment: USD -23.95'
$ws.Cells.Item(13,3).Value = '''-23.95'
$ws.Cells.Item(13,4).Value = 'PayPal'
$ws.Cells.Item(14,1).Value = '''12/14/2023'
$ws.Cells.Item(14,2).Value = 'PreApproved Payment Bill User Payment: USD -42.00'
$ws.Cells.Item(14,3).Value = '''-42.00'
$ws.Cells.Item(14,4).Value = 'PayPal'
$ws.Cells.Item(15,1).Value = '''12/04/2023'
$ws.Cells.Item(15,2).Value = 'PreApproved Payment Bill User Payment: USD -5.23'
$ws.Cells.Item(15,3).Value = '''-5.23'
$ws.Cells.Item(15,4).Value = 'PayPal'
$ws.Cells.Item(16,1).Value = '''12/19/2023'
$ws.Cells.Item(16,2).Value = 'PreApproved Payment Bill User Payment: USD -37.60'
$ws.Cells.Item(16,3).Value = '''-37.60'
$ws.Cells.Item(16,4).Value = 'PayPal'
$ws.Cells.Item(17,1).Value = '''12/01/2023'
$ws.Cells.Item(17,2).Value = 'PreApproved Payment Bill User Payment: USD -107.95'
$ws.Cells.Item(17,3).Value = '''-107.95'
$ws.Cells.Item(17,4).Value = 'PayPal'
$ws.Cells.Item(18,1).Value = '''12/09/2023'
$ws.Cells.Item(18,2).Value = 'PreApproved Payment Bill User Payment: USD -50.00'
$ws.Cells.Item(18,3).Value = '''-50.00'
$ws.Cells.Item(18,4).Value = 'PayPal'

$ws = $wb.Worksheets.Item('eBay')
$ws.Cells.Item(2,1).Value = 'Feb 17, 2024'
$ws.Cells.Item(2,2).Value = 'Delivered on Wed, Feb 21 View similar items Returns not accepted. View seller''s other items (Large Lot) Original Konami Yu-Gi-Oh! YuGiOh Collectible Tins & Card US $54.00 More actions Sold by:tanks-treasures Delivered'
$ws.Cells.Item(2,3).Value = '$70.44 '
$ws.Cells.Item(2,4).Value = 'eBay'
$ws.Cells.Item(3,1).Value = 'Apr 06, 2024'
$ws.Cells.Item(3,2).Value = 'Delivered on Tue, Apr 16 Leave feedback Returns not accepted. View seller''s other items 8,000+ POKÉMON CARD BULK LOT (cards, markers, sleeves, and more!) US $170.00 More actions *Important info Sold by:chlri_65 Delivered'
$ws.Cells.Item(3,3).Value = '$275.60 '
$ws.Cells.Item(3,4).Value = 'eBay'
$ws.Cells.Item(4,1).Value = 'Nov 20, 2023'
$ws.Cells.Item(4,2).Value = 'Delivered on Mon, Dec 4 View similar items Returns not accepted. View seller''s other items Mike Schmidt Signed 1992 Program MLB Pro- Celebrity Challenge Michael Jordan Auto US $99.00 More actions Sold by:shellyscollectibles Delivered'
$ws.Cells.Item(4,3).Value = '$114.01 '
$ws.Cells.Item(4,4).Value = 'eBay'
$ws.Cells.Item(5,1).Value = 'Aug 21, 2023'
$ws.Cells.Item(5,2).Value = 'Delivered on Sat, Aug 26 View similar items Return Window Closed on Sep 24. View seller''s other items 2021-2023 Tesla Model 3 / Y Left driver side headlight assembly OEM 1514952-00-D US $148.00 More actions Sold by:cali_autoparts Delivered'
$ws.Cells.Item(5,3).Value = '$156.88 '
$ws.Cells.Item(5,4).Value = 'eBay'
$ws.Cells.Item(6,1).Value = 'Nov 10, 2023'
$ws.Cells.Item(6,2).Value = 'Shipped: Est. delivery Tue, Nov 28 - Tue, Buy again Jan 2 Return Window Closed on Mar 2. View seller''s other items 2 X Kottakkal Vatagajankusa Rasam 30 Tablet | Relife For Joint | FREE SHIP More actions Quantity : 4 US $79.60 Sold by:ciriexport Delivered'
$ws.Cells.Item(6,3).Value = '$71.64 '
$ws.Cells.Item(6,4).Value = 'eBay'
$ws.Cells.Item(7,1).Value = 'Mar 02, 2023'
$ws.Cells.Item(7,2).Value = 'Delivered on Tue, Apr 4 Buy again Return Window Closed on May 4. View seller''s other items Ayurveda Kottakkal Arya Vaidya Sala RASNAIRANDADI KWATHAM (TABLET) - 100NOS Quantity : 4 More actions US $120.00 Sold by:creation_india Delivered'
$ws.Cells.Item(7,3).Value = '$120.00 '
$ws.Cells.Item(7,4).Value = 'eBay'
$ws.Cells.Item(8,1).Value = 'Aug 07, 2023'
$ws.Cells.Item(8,2).Value = 'https://www.ebay.com/mye/myebay/purchase Page 7 of 10'
$ws.Cells.Item(8,3).Value = '$63.55 '
$ws.Cells.Item(8,4).Value = 'eBay'
$ws.Cells.Item(9,1).Value = 'Apr 06, 2024'
$ws.Cells.Item(9,2).Value = 'Delivered on Thu, Apr 11 Leave feedback Return Window Closed on May 11. View seller''s other items Complete Original 1999 Pokemon Jungle Set! 64 Cards Most In Near Mint Condition! US $499.97 More actions Sold by:wholesale_gaming_store Delivered'
$ws.Cells.Item(9,3).Value = '$529.97 '
$ws.Cells.Item(9,4).Value = 'eBay'
$ws.Cells.Item(10,1).Value = 'Apr 17, 2023'
$ws.Cells.Item(10,2).Value = 'Delivered on Fri, Apr 28 Buy again Return Window Closed on May 28. View seller''s other items Indian Women''s Cotton Printed Night Gown Nighty Combo Pack of 8 Free Size US $68.00 More actions Sold by:homedecorhome Delivered'
$ws.Cells.Item(10,3).Value = '$73.00 '
$ws.Cells.Item(10,4).Value = 'eBay'
$ws.Cells.Item(11,1).Value = 'Mar 02, 2023'
$ws.Cells.Item(11,2).Value = 'Delivered on Tue, Apr 4 View similar items Return Window Closed on May 4. View seller''s other items 2 x Ayurveda Kottakkal Arya Vaidya Sala Yogaraja Gulgulu Vatika 100 Tablets Quantity : 2 More actions GBP 50.00 Sold by:creation_india 1 2 3 4 5 https://www.ebay.com/mye/myebay/purchase Page 9 of 10'
$ws.Cells.Item(11,3).Value = '$GBP 53.00 '
$ws.Cells.Item(11,4).Value = 'eBay'
$ws.Cells.Item(12,1).Value = 'Feb 19, 2024'
$ws.Cells.Item(12,2).Value = 'Delivered on Thu, Feb 22 View similar items Returns not accepted. View seller''s other items Pokemon Base Set Complete Basic ENERGY 6 Card Lot #97-102 Quantity : 4 More actions US $16.00 Sold by:xxwhitelionxx21 Delivered'
$ws.Cells.Item(12,3).Value = '$17.65 '
$ws.Cells.Item(12,4).Value = 'eBay'
$ws.Cells.Item(13,1).Value = 'Apr 10, 2023'
$ws.Cells.Item(13,2).Value = 'https://www.ebay.com/mye/myebay/purchase Page 8 of 10'
$ws.Cells.Item(13,3).Value = '$120.00 '
$ws.Cells.Item(13,4).Value = 'eBay'
$ws.Cells.Item(14,1).Value = 'Dec 16, 2023'
$ws.Cells.Item(14,2).Value = 'https://www.ebay.com/mye/myebay/purchase Page 3 of 10'
$ws.Cells.Item(14,3).Value = '$2,861.97 '
$ws.Cells.Item(14,4).Value = 'eBay'
$ws.Cells.Item(15,1).Value = 'Oct 15, 2023'
$ws.Cells.Item(15,2).Value = 'https://www.ebay.com/mye/myebay/purchase Page 6 of 10'
$ws.Cells.Item(15,3).Value = '$10.59 '
$ws.Cells.Item(15,4).Value = 'eBay'
$ws.Cells.Item(16,1).Value = 'Apr 09, 2024'
$ws.Cells.Item(16,2).Value = 'https://www.ebay.com/mye/myebay/purchase Page 1 of 10'
$ws.Cells.Item(16,3).Value = '$13.24 '
$ws.Cells.Item(16,4).Value = 'eBay'
$ws.Cells.Item(17,1).Value = 'Nov 10, 2023'
$ws.Cells.Item(17,2).Value = 'https://www.ebay.com/mye/myebay/purchase Page 5 of 10'
$ws.Cells.Item(17,3).Value = '$32.37 '
$ws.Cells.Item(17,4).Value = 'eBay'
$ws.Cells.Item(18,1).Value = 'Oct 24, 2023'
$ws.Cells.Item(18,2).Value = 'Delivered on Fri, Oct 27 View similar items Return Window Closed on Nov 10. View seller''s other items MICHEAL JORDAN #23 BULLS SIGNED NORTH CAROLINA BASKETBALL HOLOGRAM AUTHENTICATED More actions US $255.00 Sold by:joeycannellajr Delivered'
$ws.Cells.Item(18,3).Value = '$291.50 '
$ws.Cells.Item(18,4).Value = 'eBay'
$ws.Cells.Item(19,1).Value = 'Nov 28, 2023'
$ws.Cells.Item(19,2).Value = 'Order total:US $954.00(Auto-paid) • Order number:03-10858-84089 Delivered on Fri, Dec 1 View similar items Returns not accepted. View seller''s other items Bitmain Antminer S19J PRO 96T Asic Bitcoin Miner US $900.00 More actions Sold by:ns2blockchain Delivered'
$ws.Cells.Item(19,3).Value = '$954.00(Auto-paid) '
$ws.Cells.Item(19,4).Value = 'eBay'
$ws.Cells.Item(20,1).Value = 'Nov 11, 2023'
$ws.Cells.Item(20,2).Value = 'Order total:US $120.00(Auto-paid) • Order number:20-10768-49923 https://www.ebay.com/mye/myebay/purchase Page 4 of 10'
$ws.Cells.Item(20,3).Value = '$120.00(Auto-paid) '
$ws.Cells.Item(20,4).Value = 'eBay'
$ws.Cells.Item(21,1).Value = 'Oct 23, 2023'
$ws.Cells.Item(21,2).Value = 'Delivered on Sat, Oct 28 View similar items Returns not accepted. View seller''s other items TESLA MODEL 3/Y LEFT HEADLIGHT MATRIX LED GLOBAL DRIVER HEAD LIGHT LAMP OEM US $269.00 More actions Sold by:zapaska Delivered'
$ws.Cells.Item(21,3).Value = '$306.34 '
$ws.Cells.Item(21,4).Value = 'eBay'
$ws.Cells.Item(22,1).Value = 'Oct 09, 2023'
$ws.Cells.Item(22,2).Value = 'Order total:US $34.65(Auto-paid) • Order number:01-10643-90115 Delivered on Thu, Oct 19 View similar items Returns not accepted. View seller''s other items NEW JELLYPOP Off-White Dress Dressy Shoes Lace with side bow 1.5" Heel Women 7M US $26.00 More actions Sold by:wearetwo2003 Delivered'
$ws.Cells.Item(22,3).Value = '$34.65(Auto-paid) '
$ws.Cells.Item(22,4).Value = 'eBay'
$ws.Cells.Item(23,1).Value = 'Apr 06, 2024'
$ws.Cells.Item(23,2).Value = 'https://www.ebay.com/mye/myebay/purchase Page 2 of 10'
$ws.Cells.Item(23,3).Value = '$56.86 '
$ws.Cells.Item(23,4).Value = 'eBay'
$ws.Cells.Item(24,1).Value = 'Nov 10, 2023'
$ws.Cells.Item(24,2).Value = 'Delivered on Tue, Dec 19 Buy again Return Window Closed on Jan 18. View seller''s other items Kottakkal Yogaraja Gulgulu Vatika 100 tablets Free Shipping Quantity : 3 More actions US $51.33 Sold by:hometradeexport Delivered'
$ws.Cells.Item(24,3).Value = '$51.33 '
$ws.Cells.Item(24,4).Value = 'eBay'
$ws.Cells.Item(25,1).Value = 'Apr 13, 2024'
$ws.Cells.Item(25,2).Value = 'Sizes Selling Delivered on Wed, Apr 17 Leave feedback Returns not accepted. Collection beta View seller''s other items Pokémon bulk cards lot (Unknown Amount Swipe To The eBay vault See More Pics) US $65.00 More actions Sold by:antant1313 Delivered'
$ws.Cells.Item(25,3).Value = '$91.04 '
$ws.Cells.Item(25,4).Value = 'eBay'
$ws.Cells.Item(26,1).Value = 'Apr 12, 2023'
$ws.Cells.Item(26,2).Value = 'Delivered on Fri, May 5 Buy again Return Window Closed on Jul 4. View seller''s other items Kottakkal Vatagajankusa Rasam 250g 30 Tablet | Relife For Joint | FREE SHIPPING Quantity : 5 More actions US $50.00 Sold by:ayurvedashop Delivered'
$ws.Cells.Item(26,3).Value = '$50.00 '
$ws.Cells.Item(26,4).Value = 'eBay'
